$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Column letters for the 10 period columns (D..M)
$cols = @("D","E","F","G","H","I","J","K","L","M")

# Ordered list of (row, 10 values D..M) pairs, matching the quarterly figures
# added to the database (rows 15 and 23 move from the literal "-" placeholder
# to a numeric 0, so every cell in this block ends up a plain number).
$rowData = @(
    , @(11, @(508190, 786232, 772894, 973393, 1327612, 1880321, 1816861, 1228832, 2079427, 3032253))
    , @(12, @(-405451, -661942, -655362, -729448, -821767, -1525236, -1450527, -888961, -1741535, -2761434))
    , @(13, @(102739, 124290, 117532, 243945, 505845, 355085, 366334, 339871, 337892, 270819))
    , @(14, @(-17989, -40058, -28254, -31755, -35424, -48336, -62519, -75665, -73116, -70652))
    , @(15, @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0))
    , @(16, @(-2731, 3850, 7046, 3977, 3139, 6404, 1019, 2280, 9220, 13199))
    , @(17, @(82019, 88082, 96324, 216167, 473560, 313153, 304834, 266486, 273996, 213366))
    , @(18, @(-10859, -22968, -28006, -36859, -40612, -32909, -31674, -27108, -37120, -41094))
    , @(19, @(594592, 81989, 352845, 3530, 4514, 24270, 154685, 21246, 8611, 23336))
    , @(20, @(665752, 147103, 421163, 182838, 437462, 304514, 427845, 260624, 245487, 195608))
    , @(21, @(-45787, -25969, 21335, -41139, -98161, -65648, 204948, -43908, -66461, -38255))
    , @(22, @(619965, 121134, 442498, 141699, 339301, 238866, 632793, 216716, 179026, 157353))
    , @(23, @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0))
    , @(24, @(619965, 121134, 442498, 141699, 339301, 238866, 632793, 216716, 179026, 157353))
    , @(25, @(134, 26, 96, 31, 74, 52, 137, 47, 39, 26))
    , @(26, @(4612528, 4612528, 4612528, 4612528, 4612528, 4612528, 4612528, 4612528, 4612528, 6103562))
    , @(27, @(102, 20, 73, 23, 56, 39, 104, 36, 29, 26))
)

foreach ($pair in $rowData) {
    $r = $pair[0]
    $vals = $pair[1]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
